$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows with new manager credentials (rows 3 and 5 changed)
$ws.Range("A3").Value = "mngr83460"
$ws.Range("B3").Value = "qAbUzyj"

$ws.Range("A5").Value = "mngr164225"
$ws.Range("B5").Value = "jahetAp"

# Update the selected cell in the sheet view
$ws.Range("A6").Select()
